$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 1944
    $ws.Range("F4").Value = 1204
    $ws.Range("F5").Value = 1351
    $ws.Range("F7").Value = 6063
}
